$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The 2020 column (K) gets duplicated into a brand-new column L, carrying
# over both the value and the full cell formatting (style) from column K.
$ws.Range("K3").Copy($ws.Range("L3"))
$ws.Range("K4").Copy($ws.Range("L4"))

# Move/record the active selection at M12, as reflected in the saved
# sheet view.
$ws.Range("M12").Select()
